$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Fill "Name Ramp" column (C) for both new rows first, matching original authoring order
$ws.Cells.Item(8, 3).Value = "7_JM_MS"
$ws.Cells.Item(9, 3).Value = "8_JM_MS"

# Fill "Name" column (B) for both new rows next
$ws.Cells.Item(8, 2).Value = "221120_xgb_tuned"
$ws.Cells.Item(9, 2).Value = "221121_xgb_tuned"

# Fill "Date" column (A) with the same numeric date serials / number format as existing rows
$ws.Cells.Item(8, 1).Value2 = 44885
$ws.Cells.Item(8, 1).NumberFormat = "m/d/yy"
$ws.Cells.Item(9, 1).Value2 = 44886
$ws.Cells.Item(9, 1).NumberFormat = "m/d/yy"

# Fill "Hand in" column (D) by copying an existing "TRUE" text cell so it stays text, not boolean
$ws.Cells.Item(2, 4).Copy()
$ws.Cells.Item(8, 4).PasteSpecial(-4163)
$ws.Cells.Item(2, 4).Copy()
$ws.Cells.Item(9, 4).PasteSpecial(-4163)

# Fill "By" column (E)
$ws.Cells.Item(8, 5).Value = "Maria"
$ws.Cells.Item(9, 5).Value = "Maria"

# Resize the table to include the two new rows
$tbl = $ws.ListObjects.Item("Tabelle2")
$tbl.Resize($ws.Range("A1:E9"))

$ws.Range("B10").Select()
